# "Actualizar" bot update: appends the latest IPSA daily closing values
# (rows 760-768, 2021-01-28 through 2021-02-05) to the IPSA sheet, the same
# way the scheduled importer does every run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IPSA")

# FECHA / IPSA value for each new trading day. "--" marks a day with no
# published value (same convention already used throughout column B, e.g.
# rows 755/756).
$newRows = @(
    @{ Row = 760; Date = "01/28/2021"; Value = 4410.8 },
    @{ Row = 761; Date = "01/29/2021"; Value = 4288.6499999999996 },
    @{ Row = 762; Date = "01/30/2021"; Value = "--" },
    @{ Row = 763; Date = "01/31/2021"; Value = "--" },
    @{ Row = 764; Date = "02/01/2021"; Value = 4381.01 },
    @{ Row = 765; Date = "02/02/2021"; Value = 4416.95 },
    @{ Row = 766; Date = "02/03/2021"; Value = 4401.87 },
    @{ Row = 767; Date = "02/04/2021"; Value = 4412 },
    @{ Row = 768; Date = "02/05/2021"; Value = 4449.45 }
)

$lastOldRow = 759

foreach ($r in $newRows) {
    $rowNum = $r.Row
    $dateCell = $ws.Cells.Item($rowNum, 1)
    $valueCell = $ws.Cells.Item($rowNum, 2)

    # Copy formatting down from the previous block so the appended cells
    # reuse the existing style records instead of minting new ones: A759 for
    # the date column, B759 for numeric closes, B755 for the "--" text days.
    $ws.Range("A" + $lastOldRow).Copy() | Out-Null
    $dateCell.PasteSpecial(-4122) | Out-Null

    if ($r.Value -eq "--") {
        $ws.Range("B755").Copy() | Out-Null
    } else {
        $ws.Range("B" + $lastOldRow).Copy() | Out-Null
    }
    $valueCell.PasteSpecial(-4122) | Out-Null

    $dateCell.Value = $r.Date
    $valueCell.Value = $r.Value
}
$excel.CutCopyMode = $false

# Grow the "IPSA" defined name to cover the newly appended rows.
$wb.Names.Item("IPSA").RefersTo = "=IPSA!`$A`$1:`$B`$768"

# Move the selection down to the new last cell, matching what Excel leaves
# selected after appending rows at the bottom of the frozen-pane view.
$ws.Activate()
$ws.Range("B768").Select()
